# Update cryptocurrency price/volume data per the scraped source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.881.51"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.779.08"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.28"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4500"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3592"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07494"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.33"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.102"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.96"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.071"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.246"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.786.96"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.64"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06457"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.842"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.954.10"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.39"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.090"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.44"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.36"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.992.09"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.232"
$ws.Range("E29").Value = "  +6.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.42"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.103"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09192"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.675"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.601"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.95"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02299"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06141"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2098"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6352"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.975"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.189"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.405"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.966"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.40"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5930"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.739"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.88"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.965"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06943"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.145"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.99"
$ws.Range("E51").Value = "  +0.90%  "
